# edit.ps1
# Applies the 'Atualizacao de bases das ligas, do dia: 02-05-2024 as 20:28' update
# to the 'Colombia Primera A' sheet: rows 424-430 are re-shuffled (odds data for the
# same match-day got reassigned to different id/fixture rows) and rows 432-435 get
# small odds corrections.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Colombia Primera A")

# --- Rows 424-430: full row re-shuffle (columns B and E:AB; A/C/D stay as-is) ---
$row424 = @(7736841,'Colombia Primera A',45410.79166666666,'Atletico Bucaramanga','Alianza',1,0,'H',1.666,3.5,5,1.65,3.75,5.75,-0.75,1.8,2,2.25,1.9,1.9,0.6499999999999999,-1,-1,0.4,-0.5,-1,0.8999999999999999)
$row425 = @(7658989,'Colombia Primera A',45410.79166666666,'Jaguares de Cordoba','Independiente Santa Fe',1,0,'H',3,3.2,2.3,3.4,3.6,2.05,0.25,2,1.8,2.5,1.8,2,2.4,-1,-1,1,-1,-1,1)
$row426 = @(7658987,'Colombia Primera A',45410.79166666666,'Deportivo Cali','Junior',0,0,'D',2.7,3.25,2.4,3.2,3.1,2.4,0.25,1.8,2.05,2.25,1.975,1.875,-1,2.1,-1,0.4,-0.5,-1,0.875)
$row427 = @(7658985,'Colombia Primera A',45410.79166666666,'Aguilas Doradas','Fortaleza',1,1,'D',1.75,3.2,5,2.05,3.2,4,-0.5,2.025,1.775,2,1.8,2,-1,2.2,-1,-1,0.7749999999999999,0,0)
$row428 = @(7658915,'Colombia Primera A',45410.79166666666,'Once Caldas','America de Cali',0,0,'D',2.3,3,3.1,2.3,3.2,3.3,-0.25,1.975,1.875,2.25,2.025,1.825,-1,2.2,-1,-0.5,0.4375,-1,0.825)
$row429 = @(7658988,'Colombia Primera A',45410.79166666666,'Envigado FC','Independiente Medellin',0,1,'A',4.2,3.4,1.8,5.25,3.6,1.7,0.75,1.925,1.875,2.25,1.775,2.025,-1,-1,0.7,-0.5,0.4375,-1,1.025)
$row430 = @(7658914,'Colombia Primera A',45410.79166666666,'La Equidad','Deportivo Pereira',0,2,'A',2,3.1,3.75,2.25,3.2,3.3,-0.25,1.925,1.875,2,1.825,1.975,-1,-1,2.3,-1,0.875,0,0)

$arr = New-Object 'object[,]' 1,27
for ($i = 0; $i -lt 27; $i++) { $arr[0,$i] = $row424[$i] }
$ws.Range("B424:AB424").Value2 = $arr

$arr = New-Object 'object[,]' 1,27
for ($i = 0; $i -lt 27; $i++) { $arr[0,$i] = $row425[$i] }
$ws.Range("B425:AB425").Value2 = $arr

$arr = New-Object 'object[,]' 1,27
for ($i = 0; $i -lt 27; $i++) { $arr[0,$i] = $row426[$i] }
$ws.Range("B426:AB426").Value2 = $arr

$arr = New-Object 'object[,]' 1,27
for ($i = 0; $i -lt 27; $i++) { $arr[0,$i] = $row427[$i] }
$ws.Range("B427:AB427").Value2 = $arr

$arr = New-Object 'object[,]' 1,27
for ($i = 0; $i -lt 27; $i++) { $arr[0,$i] = $row428[$i] }
$ws.Range("B428:AB428").Value2 = $arr

$arr = New-Object 'object[,]' 1,27
for ($i = 0; $i -lt 27; $i++) { $arr[0,$i] = $row429[$i] }
$ws.Range("B429:AB429").Value2 = $arr

$arr = New-Object 'object[,]' 1,27
for ($i = 0; $i -lt 27; $i++) { $arr[0,$i] = $row430[$i] }
$ws.Range("B430:AB430").Value2 = $arr

# --- Rows 432-435: individual odds corrections ---
$ws.Range("M432").Value2 = 2.05
$ws.Range("Q432").Value2 = 1.775
$ws.Range("R432").Value2 = 2.1
$ws.Range("S432").Value2 = 2.25
$ws.Range("T432").Value2 = 2.05
$ws.Range("U432").Value2 = 1.8
$ws.Range("Q433").Value2 = 1.85
$ws.Range("R433").Value2 = 2
$ws.Range("Q434").Value2 = 2
$ws.Range("R434").Value2 = 1.85
$ws.Range("T434").Value2 = 2
$ws.Range("U434").Value2 = 1.85
$ws.Range("T435").Value2 = 1.925
$ws.Range("U435").Value2 = 1.925
